# Update the NATMI ligand-receptor edge table with newly recomputed TPM-derived
# expression values. The per-cluster average/total expression values changed
# (new TPM normalization), which cascades into the derived-specificity and
# edge-weight columns for every sending/target cluster pair row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand/receptor average expression value per cluster (column G / M source).
$newAvg = @{
    "ECs"           = 116.537051
    "FAPs"          = 19.342779
    "MuSCs"         = 3.680496333333333
    "Resolving-Mac" = 0.4968726666666667
}

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# Total expression value = average expression value * number of expressing cells (3).
$newTotal = @{}
foreach ($c in $clusters) {
    $newTotal[$c] = $newAvg[$c] * 3
}

# Derived specificity of average/total expression value = cluster value / sum across clusters.
$sumAvg = 0
$sumTotal = 0
foreach ($c in $clusters) {
    $sumAvg = $sumAvg + $newAvg[$c]
    $sumTotal = $sumTotal + $newTotal[$c]
}

$specAvg = @{}
$specTotal = @{}
foreach ($c in $clusters) {
    $specAvg[$c] = $newAvg[$c] / $sumAvg
    $specTotal[$c] = $newTotal[$c] / $sumTotal
}

# Sum of edge average/total expression weights across every sending x target pair,
# needed to derive the edge specificity columns (S, T).
$sumEdgeAvg = 0
$sumEdgeTotal = 0
foreach ($i in $clusters) {
    foreach ($j in $clusters) {
        $sumEdgeAvg = $sumEdgeAvg + ($newAvg[$i] * $newAvg[$j])
        $sumEdgeTotal = $sumEdgeTotal + ($newTotal[$i] * $newTotal[$j])
    }
}

for ($row = 2; $row -le 17; $row++) {
    $sendCluster = $ws.Cells.Item($row, 1).Text
    $targetCluster = $ws.Cells.Item($row, 4).Text

    $gVal = $newAvg[$sendCluster]
    $hVal = $newTotal[$sendCluster]
    $iVal = $specAvg[$sendCluster]
    $jVal = $specTotal[$sendCluster]

    $mVal = $newAvg[$targetCluster]
    $nVal = $newTotal[$targetCluster]
    $oVal = $specAvg[$targetCluster]
    $pVal = $specTotal[$targetCluster]

    $qVal = $gVal * $mVal
    $rVal = $hVal * $nVal
    $sVal = $qVal / $sumEdgeAvg
    $tVal = $rVal / $sumEdgeTotal

    $ws.Cells.Item($row, 7).Value = $gVal    # G: Ligand average expression value
    $ws.Cells.Item($row, 8).Value = $hVal    # H: Ligand total expression value
    $ws.Cells.Item($row, 9).Value = $iVal    # I: Ligand derived specificity of average expr
    $ws.Cells.Item($row, 10).Value = $jVal   # J: Ligand derived specificity of total expr

    $ws.Cells.Item($row, 13).Value = $mVal   # M: Receptor average expression value
    $ws.Cells.Item($row, 14).Value = $nVal   # N: Receptor total expression value
    $ws.Cells.Item($row, 15).Value = $oVal   # O: Receptor derived specificity of average expr
    $ws.Cells.Item($row, 16).Value = $pVal   # P: Receptor derived specificity of total expr

    $ws.Cells.Item($row, 17).Value = $qVal   # Q: Edge average expression weight
    $ws.Cells.Item($row, 18).Value = $rVal   # R: Edge total expression weight
    $ws.Cells.Item($row, 19).Value = $sVal   # S: Edge average expression derived specificity
    $ws.Cells.Item($row, 20).Value = $tVal   # T: Edge total expression derived specificity
}
